# Consist updates, [747] WIP
#
# Adds a new worksheet "HSL" after the existing "ZvNL" sheet, with a single
# note cell (A1) describing a bug ("Shl1/2-Asra sein op rood, tab ok"), sizes
# column A to fit, leaves the selection on A2 (as Excel does after typing
# into A1 and hitting Enter), and makes the new sheet the active tab.

$wb = $excel.ActiveWorkbook

$zvnl = $wb.Worksheets.Item(1)

# Insert the new sheet right after "ZvNL" so it becomes the 2nd tab.
$hsl = $wb.Worksheets.Add($null, $zvnl)
$hsl.Name = "HSL"

$hsl.Range("A1").Value = "Shl1/2-Asra sein op rood, tab ok"

# Widen column A to fit the note text.
$hsl.Columns.Item(1).ColumnWidth = 54.66

# Match the post-entry selection state (cursor moved down to A2) and make
# the new sheet the active one.
$hsl.Range("A2").Select() | Out-Null
$hsl.Activate()
